$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.848.00'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.564.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.92%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.563.14'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.158'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.91%  '
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("E12").Value = '  +1.27%  '
$ws.Range("E13").Value = '  +1.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.042.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.31'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.74%  '
$ws.Range("E16").Value = '  +5.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.761.45'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.564.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.68'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.64%  '
$ws.Range("E20").Value = '  +1.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '366.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.79%  '
$ws.Range("E22").Value = '  +2.02%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.32'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.76'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("E27").Value = '  +1.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.690.77'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0929'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '521.10'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.80'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.71%  '
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -1.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.12'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.98'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.28%  '
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("E41").Value = '  +1.01%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.97'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.09%  '
$ws.Range("E44").Value = '  -1.20%  '
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.09'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '153.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.18%  '
$ws.Range("E48").Value = '  +2.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.525'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("E50").Value = '  -0.71%  '
$ws.Range("E51").Value = '  +1.86%  '
